$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column values are plain numeric-looking strings (e.g. "1.006",
# "26.841.85") that must stay as literal text like the rest of the sheet -
# force Text number format first so Excel does not coerce them to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.841.85'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.829.23'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.28'
$ws.Range("E5").Value = '  -1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4617'
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3690'
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07173'
$ws.Range("E9").Value = '  -2.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8758'
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07841'
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.57'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.872.76'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.328'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.383'
$ws.Range("E15").Value = '  -2.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.30'
$ws.Range("E16").Value = '  -6.22%  '
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008719'
$ws.Range("E18").Value = '  -1.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.877.84'
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.45'
$ws.Range("E21").Value = '  -2.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.987'
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.42'
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.976'
$ws.Range("E24").Value = '  +4.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.78'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.19'
$ws.Range("E26").Value = '  -1.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.960'
$ws.Range("E27").Value = '  -5.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.57'
$ws.Range("E28").Value = '  -2.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.922'
$ws.Range("E29").Value = '  -3.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08806'
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.128'
$ws.Range("E31").Value = '  +3.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7534'
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.460'
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.130'
$ws.Range("E34").Value = '  -2.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.559'
$ws.Range("E35").Value = '  -1.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.086'
$ws.Range("E36").Value = '  +0.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01933'
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.929'
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05118'
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.897'
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4971'
$ws.Range("E41").Value = '  -3.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1595'
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.302'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4674'
$ws.Range("E44").Value = '  -3.50%  '
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.23'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("E48").Value = '  -2.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06096'
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.39'
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.39'
$ws.Range("E51").Value = '  -1.77%  '
